$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-07-13 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-14 Sunday", 2) | Out-Null
$d.Content.Find.Execute("64×81=5184", $true, $false, $false, $false, $false, $true, 1, $false, "26×77=2002", 2) | Out-Null
$d.Content.Find.Execute("19×14=266", $true, $false, $false, $false, $false, $true, 1, $false, "76×43=3268", 2) | Out-Null
$d.Content.Find.Execute("59×80=4720", $true, $false, $false, $false, $false, $true, 1, $false, "18×81=1458", 2) | Out-Null
$d.Content.Find.Execute("65×25=1625", $true, $false, $false, $false, $false, $true, 1, $false, "23×15=345", 2) | Out-Null
$d.Content.Find.Execute("17×96=1632", $true, $false, $false, $false, $false, $true, 1, $false, "41×63=2583", 2) | Out-Null
$d.Content.Find.Execute("15×36=540", $true, $false, $false, $false, $false, $true, 1, $false, "78×35=2730", 2) | Out-Null
$d.Content.Find.Execute("61×49=2989", $true, $false, $false, $false, $false, $true, 1, $false, "11×29=319", 2) | Out-Null
$d.Content.Find.Execute("52×82=4264", $true, $false, $false, $false, $false, $true, 1, $false, "82×76=6232", 2) | Out-Null
$d.Content.Find.Execute("42×30=1260", $true, $false, $false, $false, $false, $true, 1, $false, "47×52=2444", 2) | Out-Null
$d.Content.Find.Execute("20×11=220", $true, $false, $false, $false, $false, $true, 1, $false, "97×57=5529", 2) | Out-Null
$d.Content.Find.Execute("62×46=2852", $true, $false, $false, $false, $false, $true, 1, $false, "83×88=7304", 2) | Out-Null
$d.Content.Find.Execute("13×59=767", $true, $false, $false, $false, $false, $true, 1, $false, "41×15=615", 2) | Out-Null
$d.Content.Find.Execute("84×92=7728", $true, $false, $false, $false, $false, $true, 1, $false, "98×43=4214", 2) | Out-Null
$d.Content.Find.Execute("98×45=4410", $true, $false, $false, $false, $false, $true, 1, $false, "98×91=8918", 2) | Out-Null
$d.Content.Find.Execute("57×39=2223", $true, $false, $false, $false, $false, $true, 1, $false, "31×34=1054", 2) | Out-Null
$d.Content.Find.Execute("59×57=3363", $true, $false, $false, $false, $false, $true, 1, $false, "60×60=3600", 2) | Out-Null
$d.Content.Find.Execute("80×65=5200", $true, $false, $false, $false, $false, $true, 1, $false, "88×77=6776", 2) | Out-Null
$d.Content.Find.Execute("20×93=1860", $true, $false, $false, $false, $false, $true, 1, $false, "67×94=6298", 2) | Out-Null
$d.Content.Find.Execute("22×61=1342", $true, $false, $false, $false, $false, $true, 1, $false, "61×25=1525", 2) | Out-Null
$d.Content.Find.Execute("67×26=1742", $true, $false, $false, $false, $false, $true, 1, $false, "17×94=1598", 2) | Out-Null
$d.Content.Find.Execute("14×21=294", $true, $false, $false, $false, $false, $true, 1, $false, "76×64=4864", 2) | Out-Null
$d.Content.Find.Execute("80×97=7760", $true, $false, $false, $false, $false, $true, 1, $false, "97×65=6305", 2) | Out-Null
$d.Content.Find.Execute("91×51=4641", $true, $false, $false, $false, $false, $true, 1, $false, "51×66=3366", 2) | Out-Null
$d.Content.Find.Execute("94×84=7896", $true, $false, $false, $false, $false, $true, 1, $false, "80×32=2560", 2) | Out-Null
$d.Content.Find.Execute("28×46=1288", $true, $false, $false, $false, $false, $true, 1, $false, "26×33=858", 2) | Out-Null
